$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# 1. Insert a new column at G (old G "run" shifts to H)
$ws.Columns.Item(7).Insert()

# 2. Add new data rows for subject S7 / Elvira (rows 92-106)
# This also registers the new shared strings in the same order as the target diff
for ($r = 92; $r -le 106; $r++) {
  $ws.Cells.Item($r,1).Value = "S7"
  $ws.Cells.Item($r,2).Value = "Elvira"
}

# Date column (C) - copy format (style) from an existing date cell, then set the value
$ws.Cells.Item(2,3).Copy() | Out-Null
$ws.Range("C92:C106").PasteSpecial(-4122) | Out-Null
for ($r = 92; $r -le 106; $r++) {
  $ws.Cells.Item($r,3).Value = 41360
}

for ($r = 92; $r -le 106; $r++) {
  $ws.Cells.Item($r,4).Value = "2013-03-27-elvira"
}

# E (fileName), F (condition), G (frequency), H (run) per new row
$ws.Cells.Item(92,5).Value = "2013-03-27-15-10-32"
$ws.Cells.Item(92,6).Value = "hybrid-8-57Hz"
$ws.Cells.Item(92,7).Value = 8.57
$ws.Cells.Item(92,8).Value = 1
$ws.Cells.Item(93,5).Value = "2013-03-27-15-27-12"
$ws.Cells.Item(93,6).Value = "hybrid-10Hz"
$ws.Cells.Item(93,7).Value = 10
$ws.Cells.Item(93,8).Value = 1
$ws.Cells.Item(94,5).Value = "2013-03-27-15-34-15"
$ws.Cells.Item(94,6).Value = "oddball"
$ws.Cells.Item(94,7).Value = 0
$ws.Cells.Item(94,8).Value = 1
$ws.Cells.Item(95,5).Value = "2013-03-27-15-40-06"
$ws.Cells.Item(95,6).Value = "hybrid-12Hz"
$ws.Cells.Item(95,7).Value = 12
$ws.Cells.Item(95,8).Value = 1
$ws.Cells.Item(96,5).Value = "2013-03-27-15-45-54"
$ws.Cells.Item(96,6).Value = "hybrid-12Hz"
$ws.Cells.Item(96,7).Value = 12
$ws.Cells.Item(96,8).Value = 2
$ws.Cells.Item(97,5).Value = "2013-03-27-15-58-21"
$ws.Cells.Item(97,6).Value = "hybrid-8-57Hz"
$ws.Cells.Item(97,7).Value = 8.57
$ws.Cells.Item(97,8).Value = 2
$ws.Cells.Item(98,5).Value = "2013-03-27-16-04-13"
$ws.Cells.Item(98,6).Value = "hybrid-8-57Hz"
$ws.Cells.Item(98,7).Value = 8.57
$ws.Cells.Item(98,8).Value = 3
$ws.Cells.Item(99,5).Value = "2013-03-27-16-12-09"
$ws.Cells.Item(99,6).Value = "hybrid-10Hz"
$ws.Cells.Item(99,7).Value = 10
$ws.Cells.Item(99,8).Value = 2
$ws.Cells.Item(100,5).Value = "2013-03-27-16-20-10"
$ws.Cells.Item(100,6).Value = "hybrid-12Hz"
$ws.Cells.Item(100,7).Value = 12
$ws.Cells.Item(100,8).Value = 3
$ws.Cells.Item(101,5).Value = "2013-03-27-16-26-09"
$ws.Cells.Item(101,6).Value = "hybrid-15Hz"
$ws.Cells.Item(101,7).Value = 15
$ws.Cells.Item(101,8).Value = 1
$ws.Cells.Item(102,5).Value = "2013-03-27-16-44-22"
$ws.Cells.Item(102,6).Value = "hybrid-15Hz"
$ws.Cells.Item(102,7).Value = 15
$ws.Cells.Item(102,8).Value = 2
$ws.Cells.Item(103,5).Value = "2013-03-27-16-50-29"
$ws.Cells.Item(103,6).Value = "hybrid-10Hz"
$ws.Cells.Item(103,7).Value = 10
$ws.Cells.Item(103,8).Value = 3
$ws.Cells.Item(104,5).Value = "2013-03-27-16-56-34"
$ws.Cells.Item(104,6).Value = "oddball"
$ws.Cells.Item(104,7).Value = 0
$ws.Cells.Item(104,8).Value = 2
$ws.Cells.Item(105,5).Value = "2013-03-27-17-02-11"
$ws.Cells.Item(105,6).Value = "oddball"
$ws.Cells.Item(105,7).Value = 0
$ws.Cells.Item(105,8).Value = 3
$ws.Cells.Item(106,5).Value = "2013-03-27-17-08-36"
$ws.Cells.Item(106,6).Value = "hybrid-15Hz"
$ws.Cells.Item(106,7).Value = 15
$ws.Cells.Item(106,8).Value = 3

# 3. Fill in the new frequency column (G) for existing rows 2-91, based on the condition (F)
$ws.Cells.Item(2,7).Value = 10
$ws.Cells.Item(3,7).Value = 15
$ws.Cells.Item(4,7).Value = 0
$ws.Cells.Item(5,7).Value = 15
$ws.Cells.Item(6,7).Value = 12
$ws.Cells.Item(7,7).Value = 12
$ws.Cells.Item(8,7).Value = 8.57
$ws.Cells.Item(9,7).Value = 12
$ws.Cells.Item(10,7).Value = 10
$ws.Cells.Item(11,7).Value = 8.57
$ws.Cells.Item(12,7).Value = 8.57
$ws.Cells.Item(13,7).Value = 15
$ws.Cells.Item(14,7).Value = 0
$ws.Cells.Item(15,7).Value = 10
$ws.Cells.Item(16,7).Value = 0
$ws.Cells.Item(17,7).Value = 10
$ws.Cells.Item(18,7).Value = 15
$ws.Cells.Item(19,7).Value = 8.57
$ws.Cells.Item(20,7).Value = 15
$ws.Cells.Item(21,7).Value = 10
$ws.Cells.Item(22,7).Value = 0
$ws.Cells.Item(23,7).Value = 0
$ws.Cells.Item(24,7).Value = 12
$ws.Cells.Item(25,7).Value = 12
$ws.Cells.Item(26,7).Value = 0
$ws.Cells.Item(27,7).Value = 8.57
$ws.Cells.Item(28,7).Value = 12
$ws.Cells.Item(29,7).Value = 8.57
$ws.Cells.Item(30,7).Value = 15
$ws.Cells.Item(31,7).Value = 10
$ws.Cells.Item(32,7).Value = 8.57
$ws.Cells.Item(33,7).Value = 0
$ws.Cells.Item(34,7).Value = 15
$ws.Cells.Item(35,7).Value = 8.57
$ws.Cells.Item(36,7).Value = 8.57
$ws.Cells.Item(37,7).Value = 10
$ws.Cells.Item(38,7).Value = 12
$ws.Cells.Item(39,7).Value = 12
$ws.Cells.Item(40,7).Value = 0
$ws.Cells.Item(41,7).Value = 15
$ws.Cells.Item(42,7).Value = 12
$ws.Cells.Item(43,7).Value = 15
$ws.Cells.Item(44,7).Value = 10
$ws.Cells.Item(45,7).Value = 10
$ws.Cells.Item(46,7).Value = 0
$ws.Cells.Item(47,7).Value = 10
$ws.Cells.Item(48,7).Value = 0
$ws.Cells.Item(49,7).Value = 12
$ws.Cells.Item(50,7).Value = 10
$ws.Cells.Item(51,7).Value = 15
$ws.Cells.Item(52,7).Value = 10
$ws.Cells.Item(53,7).Value = 12
$ws.Cells.Item(54,7).Value = 0
$ws.Cells.Item(55,7).Value = 8.57
$ws.Cells.Item(56,7).Value = 0
$ws.Cells.Item(57,7).Value = 15
$ws.Cells.Item(58,7).Value = 15
$ws.Cells.Item(59,7).Value = 8.57
$ws.Cells.Item(60,7).Value = 12
$ws.Cells.Item(61,7).Value = 8.57
$ws.Cells.Item(62,7).Value = 15
$ws.Cells.Item(63,7).Value = 8.57
$ws.Cells.Item(64,7).Value = 15
$ws.Cells.Item(65,7).Value = 0
$ws.Cells.Item(66,7).Value = 12
$ws.Cells.Item(67,7).Value = 10
$ws.Cells.Item(68,7).Value = 10
$ws.Cells.Item(69,7).Value = 0
$ws.Cells.Item(70,7).Value = 12
$ws.Cells.Item(71,7).Value = 15
$ws.Cells.Item(72,7).Value = 0
$ws.Cells.Item(73,7).Value = 10
$ws.Cells.Item(74,7).Value = 8.57
$ws.Cells.Item(75,7).Value = 8.57
$ws.Cells.Item(76,7).Value = 12
$ws.Cells.Item(77,7).Value = 8.57
$ws.Cells.Item(78,7).Value = 10
$ws.Cells.Item(79,7).Value = 12
$ws.Cells.Item(80,7).Value = 15
$ws.Cells.Item(81,7).Value = 0
$ws.Cells.Item(82,7).Value = 0
$ws.Cells.Item(83,7).Value = 10
$ws.Cells.Item(84,7).Value = 12
$ws.Cells.Item(85,7).Value = 15
$ws.Cells.Item(86,7).Value = 8.57
$ws.Cells.Item(87,7).Value = 12
$ws.Cells.Item(88,7).Value = 15
$ws.Cells.Item(89,7).Value = 8.57
$ws.Cells.Item(90,7).Value = 0
$ws.Cells.Item(91,7).Value = 10

# 4. Header for the new column - set LAST so "frequency" is appended at the end of sharedStrings
$ws.Cells.Item(1,7).Value = "frequency"

# 5. Column width for G to match target (stored width 10, like col F bestFit look)
$ws.Columns.Item(7).ColumnWidth = 9.166666666666666

# 6. Sheet view adjustments to match the target diff
$ws.Application.ActiveWindow.ScrollRow = 49
$ws.Range("F97").Select() | Out-Null
